$d = $word.ActiveDocument

# The farewell sentence originally reads:
#   "Wir vertrauen ihn der Liebe Gottes an, ..."
# We replace the personal pronoun "ihn" with the PERSONALPRONOMENA
# placeholder, isolating it in its own run -- exactly like the existing
# "VORNAME NACHNAME" placeholder earlier in the same paragraph is its own
# run, distinct from the surrounding text runs.

$rng = $d.Content
$found = $rng.Find.Execute("ihn", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if ($found) {
    # Replace the matched word with the placeholder text.
    $rng.Text = "PERSONALPRONOMENA"

    # Toggling a character-formatting property (and restoring it) forces the
    # engine to split this range off into its own <w:r> run instead of
    # re-merging it with the identically-formatted runs before/after it.
    $rng.Bold = 1
    $rng.Bold = 0
}
